# Weekly data refresh: insert a new "Ciboulette" observation as the most
# recent record (row 177), pushing the existing rows 177-247 down to
# 178-248.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 177 - everything from old row 177
# downward shifts down by one (old row 247 becomes new row 248).
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new weekly observation.
$ws.Range("A177").Value = 4
$ws.Range("B177").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C177").Value = "Los Lagos"
$ws.Range("D177").Value = 44784
$ws.Range("E177").Value = 10
$ws.Range("F177").Value = 100112039
$ws.Range("G177").Value = "Ciboulette"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 80
$ws.Range("K177").Value = 4000
$ws.Range("L177").Value = 4000
$ws.Range("M177").Value = 4000
$ws.Range("N177").Value = "$/docena de atados"
$ws.Range("O177").Value = "Región Metropolitana"
$ws.Range("P177").Value = 1333
$ws.Range("Q177").Value = 3
$ws.Range("R177").Value = "Hortaliza"
